# Regenerate merged AHB files
# Applies the "segment header" gray styling (already used on rows 2 and 9)
# to the remaining segment-header rows (13,17,23,27,34,40,80), and clears
# the now-obsolete "AENDERUNG" marker out of column L for all affected
# detail rows within those segments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference row that already carries the correct "segment header" look:
# A:V style = s2 (plain gray), except column B which is bold (s3),
# and column L which is centered gray with no value (s4).
$formatSource = $ws.Range("A2:V2")

# Rows whose entire A:V range needs the segment-header styling applied.
$fullHeaderRows = @(13, 17, 23, 27, 34, 40, 80)

foreach ($r in $fullHeaderRows) {
    $formatSource.Copy()
    $dst = $ws.Range("A" + $r + ":V" + $r)
    $dst.PasteSpecial(-4122)
    $ws.Range("L" + $r).ClearContents()
}

# Rows where only the column L "AENDERUNG" marker must be removed
# (style becomes the plain centered gray s4, value cleared).
$lOnlyRows = @(14, 15, 16, 18, 19, 20, 21, 22, 24, 25, 26, 28, 29, 30, 31, 32, 33, 35, 36, 37, 38, 39, 41, 42, 43)

$lFormatSource = $ws.Range("L2")
foreach ($r in $lOnlyRows) {
    $lFormatSource.Copy()
    $dst = $ws.Range("L" + $r)
    $dst.PasteSpecial(-4122)
    $dst.ClearContents()
}

$excel.CutCopyMode = 0
